$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two "period count" cells that changed from "1 période" to "2 périodes" ---
$ws.Range("C55").Value = "2 périodes"
$ws.Range("C56").Value = "2 périodes"

# --- Add the new journal entry as row 58 ---
# Copy the formatting of the row above (row 57) down onto the new row so the
# date/text/period columns keep the same styles (date format, wrap text, etc.)
$null = $ws.Range("A57:C57").Copy()
$null = $ws.Range("A58:C58").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A58").Value = "3/27/2018"
$ws.Range("B58").Value = "J'ai géré la partie des quantités dans mon panier et dans la base de données. Quand j'ajoute un article dans le panier il soustrait un dans la BD. Quand je supprime un article dans le panier il ajoute la quantité que j'ai dans mon panier dans la BD. J'ai aussi amélioré le total des prix dans le panier, il va d'abord multiplier le nombre de fois que j'ai le même article dans le panier et tout additionner ensemble."
$ws.Range("C58").Value = "2 périodes"

$ws.Rows(58).RowHeight = 60

# --- Move the selection to the next empty row in column C, matching the workbook's saved cursor position ---
$null = $ws.Range("C59").Select()
